# Added "amount" and "income" (STEUERBARESEINKOMMEN) data to the fake data
# workbook, plus EGID/EWID identifier columns. This shifts the existing
# VERMOEGEN / HASEL / HASSH columns to the right to make room.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) --------------------------------------------------
# New header cells are written in the order the new shared strings need to
# be appended to the shared-string table: EGID, EWID, AMOUNT, then
# STEUERBARESEINKOMMEN (with the carried-over VERMÖGEN/HASEL/HASSH headers
# re-set along the way, reusing their existing shared-string entries).
$ws.Range("N1").Value = "EGID"
$ws.Range("O1").Value = "EWID"
$ws.Range("P1").Value = "VERMÖGEN"
$ws.Range("T1").Value = "AMOUNT"
$ws.Range("R1").Value = "HASEL"
$ws.Range("S1").Value = "HASSH"
$ws.Range("Q1").Value = "STEUERBARESEINKOMMEN"

# The brand-new header cells (Q1:T1 fall outside the sheet's original A1:P1
# extent) don't inherit the bold header style automatically, so copy it
# over from an existing header cell.
$ws.Range("M1").Copy()
$ws.Range("Q1:T1").PasteSpecial(-4122)

# --- Row 2 -----------------------------------------------------------------
$ws.Range("N2").Value = 222
$ws.Range("O2").Value = 2
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 25000
$ws.Range("R2").Value = 0
$ws.Range("S2").Value = 0
$ws.Range("T2").Value = -100

# --- Row 3 -----------------------------------------------------------------
$ws.Range("N3").Value = 1231
$ws.Range("O3").Value = 122
$ws.Range("P3").Value = 10000
$ws.Range("Q3").Value = 500000
$ws.Range("R3").Value = 1
$ws.Range("S3").Value = 0
$ws.Range("T3").Value = 777

# --- Row 4 -----------------------------------------------------------------
$ws.Range("N4").Value = 9999999
$ws.Range("O4").Value = 999
$ws.Range("P4").Value = 20000
$ws.Range("Q4").Value = 50000
$ws.Range("T4").Value = 99

# N4/O4 pick up the same (SVARotis) font formatting already used by L4/M4
# on that row.
$ws.Range("M4").Copy()
$ws.Range("N4:O4").PasteSpecial(-4122)

# Restore the active selection to match the saved workbook state.
$ws.Range("Q10").Select()
